$d = $word.ActiveDocument

$pairs = @(
    @("890×2=1780", "319×3=957"),
    @("902×3=2706", "263×9=2367"),
    @("361×3=1083", "761×7=5327"),
    @("230×4=920", "681×7=4767"),
    @("472×9=4248", "674×6=4044"),
    @("591×6=3546", "961×9=8649"),
    @("149×8=1192", "188×7=1316"),
    @("822×3=2466", "816×8=6528"),
    @("609×9=5481", "435×6=2610"),
    @("877×5=4385", "958×5=4790"),
    @("684×8=5472", "987×6=5922"),
    @("930×5=4650", "503×4=2012"),
    @("955×4=3820", "290×8=2320"),
    @("556×3=1668", "257×9=2313"),
    @("626×2=1252", "492×8=3936"),
    @("492×7=3444", "526×4=2104"),
    @("182×7=1274", "127×6=762"),
    @("261×2=522", "396×2=792"),
    @("774×6=4644", "149×5=745"),
    @("564×3=1692", "922×2=1844"),
    @("635×3=1905", "697×6=4182"),
    @("637×3=1911", "151×6=906"),
    @("326×8=2608", "356×4=1424"),
    @("882×7=6174", "416×7=2912"),
    @("644×5=3220", "797×4=3188")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
